$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark the Price/Volume columns as Text so that numeric-looking
# strings (e.g. "1.004", "3.631") are written back as text, matching the
# original inline-string cell type, instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.551.81'
$ws.Range("E2").Value = '  -2.66%  '

# Row 3
$ws.Range("D3").Value = '1.753.85'
$ws.Range("E3").Value = '  -3.23%  '

# Row 4
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").Value = '324.22'
$ws.Range("E5").Value = '  -0.80%  '

# Row 6
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.22%  '

# Row 7
$ws.Range("D7").Value = '0.4493'
$ws.Range("E7").Value = '  +3.08%  '

# Row 8
$ws.Range("D8").Value = '0.3618'
$ws.Range("E8").Value = '  -1.30%  '

# Row 9
$ws.Range("D9").Value = '0.07524'
$ws.Range("E9").Value = '  -2.02%  '

# Row 10
$ws.Range("D10").Value = '42.19'
$ws.Range("E10").Value = '  -5.86%  '

# Row 11
$ws.Range("D11").Value = '1.103'
$ws.Range("E11").Value = '  -3.43%  '

# Row 12
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.24%  '

# Row 13
$ws.Range("D13").Value = '20.74'
$ws.Range("E13").Value = '  -5.72%  '

# Row 14
$ws.Range("D14").Value = '6.049'
$ws.Range("E14").Value = '  -4.09%  '

# Row 15
$ws.Range("D15").Value = '7.221'
$ws.Range("E15").Value = '  -4.06%  '

# Row 16
$ws.Range("D16").Value = '1.756.95'
$ws.Range("E16").Value = '  -3.94%  '

# Row 17
$ws.Range("D17").Value = '92.84'
$ws.Range("E17").Value = '  -2.73%  '

# Row 18
$ws.Range("E18").Value = '  -1.38%  '

# Row 19
$ws.Range("D19").Value = '0.06427'
$ws.Range("E19").Value = '  -1.61%  '

# Row 20
$ws.Range("E20").Value = '  +0.12%  '

# Row 21
$ws.Range("D21").Value = '17.09'
$ws.Range("E21").Value = '  -1.95%  '

# Row 22
$ws.Range("D22").Value = '5.858'
$ws.Range("E22").Value = '  -6.18%  '

# Row 23
$ws.Range("D23").Value = '27.602.33'
$ws.Range("E23").Value = '  -2.54%  '

# Row 24
$ws.Range("D24").Value = '11.26'
$ws.Range("E24").Value = '  -2.75%  '

# Row 25
$ws.Range("D25").Value = '2.097'
$ws.Range("E25").Value = '  +0.73%  '

# Row 26
$ws.Range("D26").Value = '162.87'
$ws.Range("E26").Value = '  +0.43%  '

# Row 27
$ws.Range("D27").Value = '20.46'
$ws.Range("E27").Value = '  -1.48%  '

# Row 28
$ws.Range("D28").Value = '1.957.30'
$ws.Range("E28").Value = '  -3.54%  '

# Row 29
$ws.Range("D29").Value = '2.128'
$ws.Range("E29").Value = '  -6.60%  '

# Row 30
$ws.Range("D30").Value = '125.86'
$ws.Range("E30").Value = '  -2.45%  '

# Row 31
$ws.Range("D31").Value = '1.091'
$ws.Range("E31").Value = '  -9.79%  '

# Row 32
$ws.Range("D32").Value = '0.09081'
$ws.Range("E32").Value = '  -1.16%  '

# Row 33
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '3.631'
$ws.Range("E33").Value = '  +3.74%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '5.545'
$ws.Range("E34").Value = '  -6.99%  '

# Row 35
$ws.Range("D35").Value = '12.11'
$ws.Range("E35").Value = '  -6.97%  '

# Row 36
$ws.Range("D36").Value = '0.02301'
$ws.Range("E36").Value = '  -1.93%  '

# Row 37
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2103'
$ws.Range("E37").Value = '  -3.11%  '

# Row 38
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '0.6405'
$ws.Range("E38").Value = '  -2.82%  '

# Row 39
$ws.Range("D39").Value = '0.05991'
$ws.Range("E39").Value = '  -3.52%  '

# Row 40
$ws.Range("D40").Value = '4.934'
$ws.Range("E40").Value = '  -4.99%  '

# Row 41
$ws.Range("E41").Value = '  -0.09%  '

# Row 42
$ws.Range("D42").Value = '1.002'
$ws.Range("E42").Value = '  +0.22%  '

# Row 43
$ws.Range("D43").Value = '1.388'
$ws.Range("E43").Value = '  -3.09%  '

# Row 44
$ws.Range("D44").Value = '7.841'
$ws.Range("E44").Value = '  -3.68%  '

# Row 45
$ws.Range("D45").Value = '13.29'
$ws.Range("E45").Value = '  -4.13%  '

# Row 46
$ws.Range("D46").Value = '0.5899'
$ws.Range("E46").Value = '  -3.52%  '

# Row 47
$ws.Range("D47").Value = '3.712'
$ws.Range("E47").Value = '  -0.90%  '

# Row 48
$ws.Range("D48").Value = '1.969'
$ws.Range("E48").Value = '  -2.36%  '

# Row 49
$ws.Range("D49").Value = '121.76'
$ws.Range("E49").Value = '  -3.28%  '

# Row 50
$ws.Range("D50").Value = '1.159'
$ws.Range("E50").Value = '  +0.24%  '

# Row 51
$ws.Range("D51").Value = '0.06877'
$ws.Range("E51").Value = '  -1.80%  '

# Restore the default cell style so no stray number-format/style index is
# left behind on these cells (keeps the written cells visually identical to
# before, same as the source data which used the default "General" style).
$ws.Range("D2:E51").Style = "Normal"
